$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new rows (6 and 7) with their values ---
# Row 6: LIFE=">=500" TOWERS="<2" MELIANTS="<2" TIME="<60" ENEMY_COINS="<300" POSITION="0x7"
$ws.Range("A6").Value = ">=500"
$ws.Range("B6").Value = "<2"
$ws.Range("C6").Value = "<2"
$ws.Range("D6").Value = "<60"
$ws.Range("E6").Value = "<300"
$ws.Range("F6").Value = "0x7"

# Row 7: LIFE="<500" TOWERS="<2" MELIANTS=">=2" TIME=">=60" ENEMY_COINS="<300" POSITION="-4,81x4,69"
$ws.Range("A7").Value = "<500"
$ws.Range("B7").Value = "<2"
$ws.Range("C7").Value = ">=2"
$ws.Range("D7").Value = ">=60"
$ws.Range("E7").Value = "<300"
$ws.Range("F7").Value = "-4,81x4,69"

# --- Match formatting of the row above (row 5 carries the style used for new rows) ---
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Explicit row height, matching the rest of the sheet (marks rows as customHeight)
$ws.Rows.Item(6).RowHeight = 15.75
$ws.Rows.Item(7).RowHeight = 15.75

# --- Update selection to match the post-edit UI state ---
[void]$ws.Range("A8:XFD16").Select()
